# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Re-sorts the worker/period detail table (rows 16-36) so that it is grouped
# by "Periodo Mora" (ascending: 1907..1912, 2001) with the three workers in
# the same relative order within each period, and corrects Shirley Patricia
# De Oro Palacin's "Salario Basico" (column G) from 828116 to 781242 to match
# the other workers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("CC","64521394","DIANI BERRIO DE DEL TORO","1907",33125,781242),
    @("CC","1143381314","ALFONSO ALVAREZ GARCIA","1907",33125,781242),
    @("CC","33272637","SHIRLEY PATRICIA DE ORO PALACIN","1907",33125,781242),
    @("CC","64521394","DIANI BERRIO DE DEL TORO","1908",31249,781242),
    @("CC","1143381314","ALFONSO ALVAREZ GARCIA","1908",33125,781242),
    @("CC","33272637","SHIRLEY PATRICIA DE ORO PALACIN","1908",33125,781242),
    @("CC","64521394","DIANI BERRIO DE DEL TORO","1909",31249,781242),
    @("CC","1143381314","ALFONSO ALVAREZ GARCIA","1909",33125,781242),
    @("CC","33272637","SHIRLEY PATRICIA DE ORO PALACIN","1909",33125,781242),
    @("CC","64521394","DIANI BERRIO DE DEL TORO","1910",31249,781242),
    @("CC","1143381314","ALFONSO ALVAREZ GARCIA","1910",33125,781242),
    @("CC","33272637","SHIRLEY PATRICIA DE ORO PALACIN","1910",33125,781242),
    @("CC","64521394","DIANI BERRIO DE DEL TORO","1911",31249,781242),
    @("CC","1143381314","ALFONSO ALVAREZ GARCIA","1911",33125,781242),
    @("CC","33272637","SHIRLEY PATRICIA DE ORO PALACIN","1911",33125,781242),
    @("CC","64521394","DIANI BERRIO DE DEL TORO","1912",31249,781242),
    @("CC","1143381314","ALFONSO ALVAREZ GARCIA","1912",31249,781242),
    @("CC","33272637","SHIRLEY PATRICIA DE ORO PALACIN","1912",31249,781242),
    @("CC","64521394","DIANI BERRIO DE DEL TORO","2001",21874,781242),
    @("CC","1143381314","ALFONSO ALVAREZ GARCIA","2001",21874,781242),
    @("CC","33272637","SHIRLEY PATRICIA DE ORO PALACIN","2001",21874,781242)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
